$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the createFlug comment (G9): the required Abflug date-format
# hint changes from " yyyy-mm-dd HH:mm " to a concrete example timestamp
# string, keeping the bold run in the middle of the rich text and the
# trailing "sein." unbolded.
$full = 'Methode, die eine Flug anlegt. Anfangs werden die Felder flugzeugid und mahlzeitid auf die default Felder der Tabellen Flugzeug und Mahlzeit belegt (jeweils die ID=1). Die FlugID setzt sich nach dem folgenden Muster zusammen: MH %relationid%/%flugnummer% (Die Flugnummer wird bei Neuanlage um 1 erhöht). Folglich hat der erste Flug der Relation 1  die FlugID MH1/1 der zweite Flug MH1/2 usw. Der String Abflug muss im Format"Tue Apr 17 17:46:00 CEST 2018" sein.'
$ws.Range("G9").Value = $full

# bold run: '"Tue Apr 17 17:46:00 CEST 2018" ' starts right after "...Format"
$ws.Range("G9").Characters(425, 32).Font.Bold = $true
# keep the trailing "sein." run explicitly non-bold (matches original formatting)
$ws.Range("G9").Characters(457, 5).Font.Bold = $false

# --- Row 9 now needs more vertical room for the longer comment text.
$ws.Rows.Item(9).RowHeight = 75

# --- Row 4's explicit height is no longer needed; let it return to the
# sheet's standard/default height.
$ws.Rows.Item(4).EntireRow.AutoFit()

# --- Column G is widened (and no longer marked as "best fit"/autofit).
$ws.Columns.Item(7).ColumnWidth = 120.5

# --- Selection / view: move the active selection to G9 (was G13), which
# also clears the previously scrolled "topLeftCell" state.
$ws.Range("G9").Select()
